# Regenerate save_data: update column G ("K") values for rows 2-23
# (K replaces the old "Strike#" values after regen of std/mean, calc, s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 0
    6  = 3
    7  = 0
    8  = 3
    9  = 0
    10 = 3
    11 = 2
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
